# "switched left and right"
#
# The "Tabelle1" sheet lists mic positions per plane (B1/B2/.../D17).
# The B1-plane rows (2-18) had their Plane-X / Plane-Y columns (C/D)
# mistakenly labeled "right" and the D2-plane rows (36-52) mistakenly
# labeled "left" - this swaps them back: B1 rows -> "left", D2 rows -> "right".
# (Row 18 / row 52's D column is a different plane value ("top") and is
# left untouched, matching the original data.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# B1 plane block: C2:C18 + D2:D17 were "right" -> should be "left"
$ws.Range("C2:C18").Value = "left"
$ws.Range("D2:D17").Value = "left"

# D2 plane block: C36:C52 + D36:D51 were "left" -> should be "right"
$ws.Range("C36:C52").Value = "right"
$ws.Range("D36:D51").Value = "right"

# Update the saved view state (scroll position + active selection) to
# match what the workbook was left showing after the edit.
[void]$ws.Activate()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("G43").Select()
